# Update cryptos list (price + 1h volume-change columns), matching the
# "Updated cryptos list ... with GitHub Actions" refresh.
#
# Column D (Price) holds text like "29.902.97" / "1.001" (dotted text,
# not a real number), so every D write goes through a temporary "@"
# (text) number format + Value2, then ClearFormats() to drop the
# temporary formatting again so the cell's style index is left exactly
# as it was (avoids Excel auto-converting values such as "1.001" or
# "24.70" into numbers, and avoids leaving a stray style behind).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value2 = '29.908.49'
$c.ClearFormats()
$ws.Range("E2").Value2 = '  +0.03%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value2 = '1.875.14'
$c.ClearFormats()
$ws.Range("E3").Value2 = '  -0.75%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value2 = '1.001'
$c.ClearFormats()
$ws.Range("E4").Value2 = '  +0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = '0.7409'
$c.ClearFormats()
$ws.Range("E5").Value2 = '  -2.94%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = '242.54'
$c.ClearFormats()
$ws.Range("E6").Value2 = '  -0.08%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = '1.000'
$c.ClearFormats()
$ws.Range("E7").Value2 = '  -0.06%  '

$ws.Range("E8").Value2 = '  +0.73%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value2 = '0.07219'
$c.ClearFormats()
$ws.Range("E9").Value2 = '  +0.84%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = '24.70'
$c.ClearFormats()
$ws.Range("E10").Value2 = '  -3.77%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value2 = '0.08330'
$c.ClearFormats()
$ws.Range("E11").Value2 = '  -2.13%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value2 = '0.7505'
$c.ClearFormats()
$ws.Range("E12").Value2 = '  -1.55%  '

$ws.Range("B13").Value2 = 'WrappedEther'
$ws.Range("C13").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value2 = '1.887.75'
$c.ClearFormats()
$ws.Range("E13").Value2 = '  -0.26%  '

$ws.Range("B14").Value2 = 'Polkadot'
$ws.Range("C14").Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = '5.386'
$c.ClearFormats()
$ws.Range("E14").Value2 = '  +0.28%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = '92.24'
$c.ClearFormats()
$ws.Range("E15").Value2 = '  -1.71%  '

$ws.Range("B16").Value2 = 'Uniswap'
$ws.Range("C16").Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value2 = '6.109'
$c.ClearFormats()
$ws.Range("E16").Value2 = '  -0.44%  '

$ws.Range("B17").Value2 = 'WrappedBTC'
$ws.Range("C17").Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value2 = '29.886.74'
$c.ClearFormats()
$ws.Range("E17").Value2 = '  +0.35%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value2 = '247.52'
$c.ClearFormats()
$ws.Range("E18").Value2 = '  +1.48%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = '13.56'
$c.ClearFormats()
$ws.Range("E19").Value2 = '  -1.54%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = '0.000007843'
$c.ClearFormats()
$ws.Range("E20").Value2 = '  +0.58%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value2 = '0.9992'
$c.ClearFormats()
$ws.Range("E21").Value2 = '  -0.09%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value2 = '2.130.33'
$c.ClearFormats()
$ws.Range("E22").Value2 = '  +1.03%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = '7.999'
$c.ClearFormats()
$ws.Range("E23").Value2 = '  -0.10%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = '1.001'
$c.ClearFormats()
$ws.Range("E24").Value2 = '  +0.03%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = '0.1548'
$c.ClearFormats()
$ws.Range("E25").Value2 = '  -4.26%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value2 = '9.301'
$c.ClearFormats()
$ws.Range("E26").Value2 = '  -1.11%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = '165.96'
$c.ClearFormats()
$ws.Range("E27").Value2 = '  +2.47%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value2 = '18.65'
$c.ClearFormats()
$ws.Range("E28").Value2 = '  -0.68%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value2 = '2.024'
$c.ClearFormats()
$ws.Range("E29").Value2 = '  -0.56%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value2 = '1.494'
$c.ClearFormats()
$ws.Range("E30").Value2 = '  +0.68%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = '4.575'
$c.ClearFormats()
$ws.Range("E31").Value2 = '  +1.98%  '

$ws.Range("E32").Value2 = '  +0.26%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value2 = '4.224'
$c.ClearFormats()
$ws.Range("E33").Value2 = '  +3.20%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = '0.05338'
$c.ClearFormats()
$ws.Range("E34").Value2 = '  -1.99%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value2 = '1.241'
$c.ClearFormats()
$ws.Range("E35").Value2 = '  -0.23%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value2 = '0.7501'
$c.ClearFormats()
$ws.Range("E36").Value2 = '  +1.00%  '

$ws.Range("E37").Value2 = '  +0.15%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value2 = '2.699'
$c.ClearFormats()
$ws.Range("E38").Value2 = '  +0.07%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value2 = '0.01963'
$c.ClearFormats()
$ws.Range("E39").Value2 = '  +0.92%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = '2.753'
$c.ClearFormats()
$ws.Range("E40").Value2 = '  -0.95%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = '0.4545'
$c.ClearFormats()
$ws.Range("E41").Value2 = '  +1.81%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value2 = '1.120.53'
$c.ClearFormats()
$ws.Range("E42").Value2 = '  +1.86%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = '6.129'
$c.ClearFormats()
$ws.Range("E43").Value2 = '  +1.01%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = '72.49'
$c.ClearFormats()
$ws.Range("E44").Value2 = '  -0.67%  '

$ws.Range("E45").Value2 = '  +1.55%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value2 = '104.44'
$c.ClearFormats()
$ws.Range("E46").Value2 = '  +1.52%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value2 = '1.863'
$c.ClearFormats()
$ws.Range("E48").Value2 = '  -0.32%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = '7.627'
$c.ClearFormats()
$ws.Range("E49").Value2 = '  -0.09%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = '9.491'
$c.ClearFormats()
$ws.Range("E50").Value2 = '  -2.39%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value2 = '2.020.35'
$c.ClearFormats()
$ws.Range("E51").Value2 = '  +0.83%  '
